# "Generate Report for Handback"
#
# The handback report workbook tracks, per target-locale sheet (zh-cn /
# de-de), the handoff/handback status of each localized file. This run
# records that handback/generation completed:
#   - Status moves from "Ready for handoff" to
#     "Handed back: in sync with en-US" (also reflected on the Overview
#     rollup sheet, which mirrors the same status text per language).
#   - The "Latest Handback DateTime" timestamps advance to the
#     handback-generation time.
#   - The stale "handback file is not the latest" error on the second
#     (e654c024...) row of each locale sheet is cleared now that the
#     handback is in sync.
#   - A couple of columns are widened/narrowed to fit the new text.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("K2").Value = "2016-10-21 01:29:34"

$zhcn.Range("C3").Value = $newStatus
$zhcn.Range("K3").Value = "2016-10-21 01:29:34"
$zhcn.Range("P3").Value = ""

$zhcn.Columns.Item(3).ColumnWidth = 29.17
$zhcn.Columns.Item(16).ColumnWidth = 12.84

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = $newStatus
$dede.Range("K2").Value = "2016-10-21 01:29:52"

$dede.Range("C3").Value = $newStatus
$dede.Range("K3").Value = "2016-10-21 01:29:52"
$dede.Range("P3").Value = ""

$dede.Columns.Item(3).ColumnWidth = 29.17
$dede.Columns.Item(16).ColumnWidth = 12.84

# ---------------------------------------------------------------------
# Overview sheet (rollup columns mirror the per-locale Status text)
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

$overview.Columns.Item(5).ColumnWidth = 29.17
$overview.Columns.Item(6).ColumnWidth = 29.17
